$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new "DocNumPerTopic" data row right after the header label
#        (row 5), pushing the existing 9 data rows (6-14) down to (7-15).
$ws.Rows(6).Insert()
$srcDocNum = $ws.Range("A7:J7")
$dstDocNum = $ws.Range("A6:J6")
$srcDocNum.Copy()
$dstDocNum.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$docNumRow = @(24, 6, 5, 17, 30, 13, 10, 7, 5, 23)
for ($i = 0; $i -lt $docNumRow.Length; $i++) {
    $ws.Cells.Item(6, $i + 1).Value2 = $docNumRow[$i]
}

# --- 2) Insert a new "Density" data row right after the "Density" label
#        (now row 17, having been pushed down by the insert above), pushing
#        the existing 9 density rows down by one more row.
$ws.Rows(18).Insert()
$srcDensity = $ws.Range("A19:J19")
$dstDensity = $ws.Range("A18:J18")
$srcDensity.Copy()
$dstDensity.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$densityRow = @(0.505390963203463, 0.41127232142857101, 0.43736810064935, 0.510545183982684, 0.50582724567099502, 0.49897186147186101, 0.432359307359307, 0.49975311147186102, 0.46702178030302999, 0.54040178571428499)
for ($i = 0; $i -lt $densityRow.Length; $i++) {
    $ws.Cells.Item(18, $i + 1).Value2 = $densityRow[$i]
}

# --- 3) Move the active selection to L9 (matches the post-edit cursor spot).
$ws.Range("L9").Select()

# --- 4) Page setup: portrait, paper size 9 (A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait
